# Refresh the crypto price / Volume(1h) figures to the latest scrape.
# Numeric-looking price strings are entered with a leading apostrophe so
# Excel keeps them as literal text (matching the sheet's existing
# plain-text price formatting) instead of auto-converting to a number.
# Rows 12 and 13 also swap coin identity (WrappedEther <-> Polygon
# re-ranked), so columns B-E are rewritten for both rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.881.35"
$ws.Range("E2").Value = "  -1.12%  "
$ws.Range("D3").Value = "1.895.47"
$ws.Range("E3").Value = "  -0.75%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'0.7543"
$ws.Range("E5").Value = "  +2.05%  "
$ws.Range("D6").Value = "'239.94"
$ws.Range("E6").Value = "  -1.63%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "'0.3038"
$ws.Range("E8").Value = "  -3.01%  "
$ws.Range("D9").Value = "'25.39"
$ws.Range("E9").Value = "  -5.83%  "
$ws.Range("D10").Value = "'0.06830"
$ws.Range("E10").Value = "  -1.93%  "
$ws.Range("D11").Value = "'0.07973"
$ws.Range("E11").Value = "  -0.08%  "
$ws.Range("B12").Value = "Polygon"
$ws.Range("C12").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D12").Value = "'0.7453"
$ws.Range("E12").Value = "  -4.32%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.902.95"
$ws.Range("E13").Value = "  -1.32%  "
$ws.Range("D14").Value = "'5.190"
$ws.Range("E14").Value = "  -1.75%  "
$ws.Range("D15").Value = "'91.08"
$ws.Range("E15").Value = "  -0.70%  "
$ws.Range("D16").Value = "29.892.03"
$ws.Range("E16").Value = "  -1.00%  "
$ws.Range("D17").Value = "'13.91"
$ws.Range("E17").Value = "  -2.54%  "
$ws.Range("D18").Value = "'5.933"
$ws.Range("E18").Value = "  +1.19%  "
$ws.Range("D19").Value = "'243.38"
$ws.Range("E19").Value = "  -0.35%  "
$ws.Range("D20").Value = "'0.000007716"
$ws.Range("E20").Value = "  -1.66%  "
$ws.Range("E21").Value = "  -0.06%  "
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("D23").Value = "'6.947"
$ws.Range("E23").Value = "  +4.50%  "
$ws.Range("D24").Value = "'9.217"
$ws.Range("E24").Value = "  -2.21%  "
$ws.Range("D25").Value = "'165.45"
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("D26").Value = "'18.73"
$ws.Range("E26").Value = "  -1.36%  "
$ws.Range("D27").Value = "'0.1295"
$ws.Range("E27").Value = "  +1.95%  "
$ws.Range("D28").Value = "'2.028"
$ws.Range("E28").Value = "  -3.97%  "
$ws.Range("D29").Value = "'1.397"
$ws.Range("E29").Value = "  +3.68%  "
$ws.Range("D30").Value = "'1.513"
$ws.Range("D31").Value = "'4.275"
$ws.Range("E31").Value = "  -0.90%  "
$ws.Range("D32").Value = "'4.023"
$ws.Range("E32").Value = "  -1.52%  "
$ws.Range("D33").Value = "'0.05358"
$ws.Range("E33").Value = "  +3.55%  "
$ws.Range("D34").Value = "'1.248"
$ws.Range("E34").Value = "  -3.60%  "
$ws.Range("D35").Value = "'0.7248"
$ws.Range("E35").Value = "  -2.79%  "
$ws.Range("D36").Value = "'2.717"
$ws.Range("E36").Value = "  -1.64%  "
$ws.Range("D37").Value = "'0.01910"
$ws.Range("E37").Value = "  -1.60%  "
$ws.Range("D38").Value = "'2.787"
$ws.Range("E38").Value = "  -0.23%  "
$ws.Range("D39").Value = "'6.169"
$ws.Range("E39").Value = "  -3.31%  "
$ws.Range("D40").Value = "'0.4400"
$ws.Range("E40").Value = "  -1.78%  "
$ws.Range("D41").Value = "'72.20"
$ws.Range("E41").Value = "  -4.20%  "
$ws.Range("D42").Value = "'1.000"
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("D43").Value = "'1.900"
$ws.Range("E43").Value = "  -1.89%  "
$ws.Range("D44").Value = "'0.8238"
$ws.Range("E44").Value = "  -1.46%  "
$ws.Range("D45").Value = "'100.89"
$ws.Range("E45").Value = "  -0.47%  "
$ws.Range("D46").Value = "'7.546"
$ws.Range("E46").Value = "  -1.61%  "
$ws.Range("D47").Value = "'9.730"
$ws.Range("E47").Value = "  -1.20%  "
$ws.Range("D48").Value = "2.058.15"
$ws.Range("E48").Value = "  -0.84%  "
$ws.Range("D49").Value = "'36.23"
$ws.Range("E49").Value = "  -3.42%  "
$ws.Range("D50").Value = "'0.05965"
$ws.Range("E50").Value = "  -0.52%  "
$ws.Range("D51").Value = "'1.465"
$ws.Range("E51").Value = "  +0.17%  "
